# ADD results from server
#
# The header row gains two new technology columns ("gb" right after "eb",
# and "btes" right before "ites") while the "gt" / "dgt" columns are
# dropped, so every label from "hp" onward shifts one column to the
# left/right as shown in the diff. Row 2 on every year-sheet is replaced
# with the freshly computed investment-cost figures from the server.

$wb = $excel.ActiveWorkbook

# New header row, valid for every sheet (A1..O1)
$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")

# New row-2 numeric results, per sheet (keyed by sheet/tab name), columns A..O
$rowValues = @{
    "2025" = @(39063.99109145206, 0, 483537.6274462014, 0, 2897240.114301849, 94331.34471502228, 0, 25342.77928792104, 0, 0, 0, 0, 0, 23638.06126801545, 19940.13531829346)
    "2030" = @(30846.52922536713, 0, 1495599.874611417, 0, 0, 70193.79982138964, 0, 56602.42752520426, 0, 0, 0, 0, 0, 51649.16401227913, 42574.77934331147)
    "2035" = @(242452.4252219552, 0, 943335.270081223, 0, 0, 1425.925979620855, 0, 39373.98526588717, 0, 0, 0, 0, 0, 53308.16490721726, 30023.09380555204)
    "2040" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 11578.49752443177, 0)
    "2045" = @(76705.58894163162, 1930.947398408091, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 28147.3462746636, 8312.661449003012)
    "2050" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($sheetName in @("2025", "2030", "2035", "2040", "2045", "2050")) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }

    $values = $rowValues[$sheetName]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item(2, $i + 1).Value = $values[$i]
    }
}
